$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = "Record"
$ws.Range("B68").Value = "Balanço Geral"
$ws.Range("C68").Value = "Obras"
$ws.Range("D68").Value = "2025-04-07T13:17"
$ws.Range("E68").Value = "Negativo"
$ws.Range("F68").Value = "Moradores da comunidade Tira-Gosto sofrem com vazamento de esgoto. Segundo moradores, chegaram a ficar 22 horas sem água. Exibido vídeo com morador faz apelo a Águas do Paraíba. Problema também seria porque uma obra da prefeitura não teria sido concluída. Em nota, Secretaria de Obras disse que, no decorrer das obras, foram necessários outros procedimentos que não estavam previstos inicialmente. Quanto à água, um técnico da concessionária vai ao local. "

$ws.Range("A69").Value = "Record"
$ws.Range("B69").Value = "Balanço Geral"
$ws.Range("C69").Value = "Obras"
$ws.Range("D69").Value = "2025-04-07T13:41"
$ws.Range("E69").Value = "Neutro"
$ws.Range("F69").Value = "Parte dos prédios históricos da cidade ainda precisam de manutenção. Entrevista com chefe do Inepac NF, Geovani Laurindo, e com a historiadora Sylvia Paes. Necessidade de 434 prédios tombados em Campos, sendo a maioria na área central. Semana passada, Defesa Civil começou a demolir prédio que estava ameaçado. Proprietários e poder público precisam cuidar. Representante do Inepac falou que Campos oferece desconto do IPTU. *matéria*"
